# Scheduled-runner data refresh: update market-price / profit columns
# (H: currentAveragePrice, I: currentAveragePriceNQ, J: currentAveragePriceHQ,
#  K: LevePriceNQ, L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 309.1
$ws.Range("I33").Value = 270.14285
$ws.Range("K33").Value = 270.14285
$ws.Range("M33").Value = -41.14285000000001

$ws.Range("H40").Value = 2394.6667
$ws.Range("I40").Value = 2011.1111
$ws.Range("K40").Value = 2011.1111
$ws.Range("M40").Value = -1836.1111

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H80").Value = 430.8
$ws.Range("I80").Value = 465
$ws.Range("J80").Value = 408
$ws.Range("K80").Value = 1395
$ws.Range("L80").Value = 1224
$ws.Range("M80").Value = -397
$ws.Range("N80").Value = -3220

$ws.Range("H83").Value = 430.8
$ws.Range("I83").Value = 465
$ws.Range("J83").Value = 408
$ws.Range("K83").Value = 4185
$ws.Range("L83").Value = 3672
$ws.Range("M83").Value = 807
$ws.Range("N83").Value = -13656

$ws.Range("H88").Value = 1833
$ws.Range("J88").Value = 4000
$ws.Range("L88").Value = 4000
$ws.Range("N88").Value = -4812

$ws.Range("H91").Value = 1833
$ws.Range("J91").Value = 4000
$ws.Range("L91").Value = 4000
$ws.Range("N91").Value = -6808

$ws.Range("H98").Value = 2943.2727
$ws.Range("I98").Value = 2337
$ws.Range("K98").Value = 2337
$ws.Range("M98").Value = -839

$ws.Range("H103").Value = 3043.2856
$ws.Range("I103").Value = 1575.75
$ws.Range("K103").Value = 4727.25
$ws.Range("M103").Value = -4141.25

$ws.Range("H122").Value = 2943.2727
$ws.Range("I122").Value = 2337
$ws.Range("K122").Value = 7011
$ws.Range("M122").Value = -4561

$ws.Range("H127").Value = 8595.375
$ws.Range("I127").Value = 8595.375
$ws.Range("K127").Value = 25786.125
$ws.Range("M127").Value = -20826.125

$ws.Range("H132").Value = 888.43475
$ws.Range("I132").Value = 888.43475
$ws.Range("K132").Value = 2665.30425
$ws.Range("M132").Value = -135.3042500000001

$ws.Range("H137").Value = 1486.4166
$ws.Range("I137").Value = 1362.3684
$ws.Range("J137").Value = 1957.8
$ws.Range("K137").Value = 4087.1052
$ws.Range("L137").Value = 5873.4
$ws.Range("M137").Value = -1537.1052
$ws.Range("N137").Value = -10973.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5963.091
$ws.Range("I32").Value = 3048.925
$ws.Range("K32").Value = 3048.925
$ws.Range("M32").Value = -2761.925

$ws.Range("H122").Value = 15437.167
$ws.Range("I122").Value = 7682.6
$ws.Range("J122").Value = 54210
$ws.Range("K122").Value = 23047.8
$ws.Range("L122").Value = 162630
$ws.Range("M122").Value = -20597.8
$ws.Range("N122").Value = -167530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2909.8333
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 1819.6666
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 1819.6666
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -2313.6666

$ws.Range("H80").Value = 283
$ws.Range("J80").Value = 311
$ws.Range("L80").Value = 311
$ws.Range("N80").Value = -2307

$ws.Range("H83").Value = 283
$ws.Range("J83").Value = 311
$ws.Range("L83").Value = 1555
$ws.Range("N83").Value = -11539

$ws.Range("H94").Value = 1083.129
$ws.Range("I94").Value = 637.7308
$ws.Range("K94").Value = 637.7308
$ws.Range("M94").Value = -186.7308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -765

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1382.6666
$ws.Range("J34").Value = 2203.8
$ws.Range("L34").Value = 6611.400000000001
$ws.Range("N34").Value = -6779.400000000001

$ws.Range("H39").Value = 1055.5
$ws.Range("J39").Value = 2997
$ws.Range("L39").Value = 8991
$ws.Range("N39").Value = -9579

$ws.Range("H109").Value = 988.2353000000001
$ws.Range("I109").Value = 988.2353000000001
$ws.Range("K109").Value = 2964.7059
$ws.Range("M109").Value = -1924.7059

$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 9036
$ws.Range("L122").Value = 17955
$ws.Range("M122").Value = -6586
$ws.Range("N122").Value = -22855

$ws.Range("H129").Value = 623
$ws.Range("I129").Value = 559.5
$ws.Range("K129").Value = 1678.5
$ws.Range("M129").Value = 3321.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2977.2856
$ws.Range("J107").Value = 3847.25
$ws.Range("L107").Value = 3847.25
$ws.Range("N107").Value = -7687.25

$ws.Range("H123").Value = 44509.09
$ws.Range("J123").Value = 45470
$ws.Range("L123").Value = 45470
$ws.Range("N123").Value = -50370

$ws.Range("H124").Value = 74969
$ws.Range("J124").Value = 74969
$ws.Range("L124").Value = 74969
$ws.Range("N124").Value = -84789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 15000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 15000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 15000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -15340

$ws.Range("H46").Value = 94454.55
$ws.Range("I46").Value = 3166.8333
$ws.Range("J46").Value = 203999.8
$ws.Range("K46").Value = 3166.8333
$ws.Range("L46").Value = 203999.8
$ws.Range("M46").Value = -2978.8333
$ws.Range("N46").Value = -204375.8

$ws.Range("H61").Value = 10102983
$ws.Range("I61").Value = 13890502
$ws.Range("K61").Value = 13890502
$ws.Range("M61").Value = -13890300

$ws.Range("H82").Value = 1424.25
$ws.Range("I82").Value = 870.5714
$ws.Range("K82").Value = 870.5714
$ws.Range("M82").Value = -509.5714

$ws.Range("H85").Value = 1424.25
$ws.Range("I85").Value = 870.5714
$ws.Range("K85").Value = 870.5714
$ws.Range("M85").Value = 377.4286

$ws.Range("H100").Value = 753
$ws.Range("I100").Value = 753
$ws.Range("K100").Value = 753
$ws.Range("M100").Value = -212

$ws.Range("H113").Value = 10102983
$ws.Range("I113").Value = 13890502
$ws.Range("K113").Value = 13890502
$ws.Range("M113").Value = -13888332

$ws.Range("H122").Value = 1735.625
$ws.Range("I122").Value = 1735.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5206.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2756.875
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

$ws.Range("H64").Value = 42199.8
$ws.Range("J64").Value = 42199.8
$ws.Range("L64").Value = 42199.8
$ws.Range("N64").Value = -42695.8

$ws.Range("H67").Value = 42199.8
$ws.Range("J67").Value = 42199.8
$ws.Range("L67").Value = 42199.8
$ws.Range("N67").Value = -43915.8

$ws.Range("H122").Value = 1746.8667
$ws.Range("I122").Value = 1733.6666
$ws.Range("K122").Value = 5200.9998
$ws.Range("M122").Value = -2750.9998
